$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.018908500671387
$ws.Range("B1").Value = 4.035956859588623
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.076705932617188
$ws.Range("E1").Value = 2.663103103637695
